$wb = $excel.ActiveWorkbook

# --- All data sheets share the same "header placeholder" layout: an empty
# styled cell at B1 and an empty styled cell at A2. These are removed on
# every sheet (simple model clean-up for testing).
foreach ($name in @("D", "W", "R", "O", "S")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("B1").Clear()
    $ws.Range("A2").Clear()
}

# --- Sheet "D": turn the row-label column into a numeric "year" index and
# label the corner cell "year\skill" (highlighted in yellow).
$d = $wb.Worksheets.Item("D")

$d.Range("B2").Value = "year\skill"
$d.Range("B2").Interior.Color = 65535
$d.Range("B2").HorizontalAlignment = -4108
$d.Range("B2").VerticalAlignment = -4108

$d.Range("B3").Value = 0
$d.Range("B4").Value = 1
$d.Range("B5").Value = 2
$d.Range("B6").Value = 3

# --- Selection / active-sheet bookkeeping matches the saved view state:
# sheet "D" selection moves to B2 and sheet "W" becomes the active tab.
$d.Range("B2").Select()
$wb.Worksheets.Item("W").Activate()
